# Apply the "AbilityConfigHitOnCollision" update to the enemy characters sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ability-config-key values for the two enemy rows (row 4 = actor_enemy_100, row 5 = actor_enemy_101)
$ws.Range("K4").Value = "Level:ActorConfigs:AbilityConfigHitOnCollision"
$ws.Range("K5").Value = "Level:ActorConfigs:AbilityConfigHitOnCollision"

# Widen column K so the new, longer values are fully visible
$ws.Columns.Item(11).ColumnWidth = 39.7

# Move the active selection to K8, matching where the editor left off
$ws.Range("K8").Select() | Out-Null
